$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2-90)
# from serial date 45204 to 45205 (i.e. advance by one day).
$ws.Range("C2:C90").Value = 45205
